$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.371.19"
$ws.Range("E2").Value = "  -2.99%  "
$ws.Range("D3").Value = "'3.682.95"
$ws.Range("E3").Value = "  -3.75%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'685.08"
$ws.Range("E5").Value = "  -2.63%  "
$ws.Range("D6").Value = "'160.94"
$ws.Range("E6").Value = "  -6.31%  "
$ws.Range("D7").Value = "'3.681.56"
$ws.Range("E7").Value = "  -3.82%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").Value = "'0.494"
$ws.Range("E9").Value = "  -6.17%  "
$ws.Range("D10").Value = "'0.146"
$ws.Range("E10").Value = "  -9.18%  "
$ws.Range("D11").Value = "'7.32"
$ws.Range("E11").Value = "  -2.46%  "
$ws.Range("D12").Value = "'0.436"
$ws.Range("E12").Value = "  -10.01%  "
$ws.Range("D13").Value = "'0.0000234"
$ws.Range("E13").Value = "  -7.57%  "
$ws.Range("D14").Value = "'4.302.78"
$ws.Range("E14").Value = "  -3.70%  "
$ws.Range("D15").Value = "'32.81"
$ws.Range("E15").Value = "  -9.47%  "
$ws.Range("D16").Value = "'3.680.27"
$ws.Range("E16").Value = "  -4.04%  "
$ws.Range("D17").Value = "'69.381.54"
$ws.Range("E17").Value = "  -3.03%  "
$ws.Range("E18").Value = "  -1.43%  "
$ws.Range("D19").Value = "'15.94"
$ws.Range("E19").Value = "  -9.31%  "
$ws.Range("D20").Value = "'6.47"
$ws.Range("E20").Value = "  -10.73%  "
$ws.Range("D21").Value = "'472.42"
$ws.Range("E21").Value = "  -8.55%  "
$ws.Range("D22").Value = "'9.90"
$ws.Range("E22").Value = "  -5.95%  "
$ws.Range("D23").Value = "'0.653"
$ws.Range("E23").Value = "  -9.31%  "
$ws.Range("D24").Value = "'79.68"
$ws.Range("E24").Value = "  -5.43%  "
$ws.Range("D25").Value = "'3.829.05"
$ws.Range("E25").Value = "  -3.56%  "
$ws.Range("D26").Value = "'0.0000128"
$ws.Range("E26").Value = "  -10.57%  "
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("D28").Value = "'11.15"
$ws.Range("E28").Value = "  -13.92%  "
$ws.Range("D29").Value = "'9.16"
$ws.Range("E29").Value = "  -11.93%  "
$ws.Range("E30").Value = "  -12.66%  "
$ws.Range("E31").Value = "  -11.33%  "
$ws.Range("D32").Value = "'6.71"
$ws.Range("E32").Value = "  -9.24%  "
$ws.Range("D33").Value = "'2.03"
$ws.Range("E33").Value = "  -9.52%  "
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("D35").Value = "'26.69"
$ws.Range("E35").Value = "  -9.04%  "
$ws.Range("E36").Value = "  -5.51%  "
$ws.Range("D37").Value = "'3.650.34"
$ws.Range("E37").Value = "  -3.60%  "
$ws.Range("D38").Value = "'8.25"
$ws.Range("E38").Value = "  -10.85%  "
$ws.Range("D39").Value = "'6.18"
$ws.Range("E39").Value = "  -3.18%  "
$ws.Range("E40").Value = "  -7.52%  "
$ws.Range("D42").Value = "'0.0910"
$ws.Range("E42").Value = "  -10.34%  "
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("D44").Value = "'0.948"
$ws.Range("E44").Value = "  -7.21%  "
$ws.Range("D45").Value = "'163.69"
$ws.Range("E45").Value = "  -5.52%  "
$ws.Range("D46").Value = "'48.36"
$ws.Range("E46").Value = "  -3.15%  "
$ws.Range("D47").Value = "'30.10"
$ws.Range("D48").Value = "'2.71"
$ws.Range("E48").Value = "  -17.44%  "
$ws.Range("D49").Value = "'1.31"
$ws.Range("E49").Value = "  -5.09%  "
$ws.Range("D50").Value = "'0.000277"
$ws.Range("E50").Value = "  -10.15%  "
$ws.Range("D51").Value = "'1.10"
$ws.Range("E51").Value = "  -5.25%  "
